$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 477: 世建发票1.pdf (page 1)
$ws.Cells.Item(477, 1).Value = @'
世建发票1.pdf
'@
$ws.Cells.Item(477, 2).Value = 1
$ws.Cells.Item(477, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90827520",
    "Invoice Date": "28 Aug 2023",
    "Currency": "USD",
    "Amount": 12474.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(477, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90827520",
    "Invoice Date": "28 Aug 2023",
    "Currency": "USD",
    "Amount": 12474.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(477, 5).Value = @'
[]
'@
$ws.Cells.Item(477, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世世世世(香香)有有有有 InvoiceNo. :90827520
Currency :USD
IssuedBy :VickyYan
Salesman :TONYGUO
Date :28Aug2023
Customer Code :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK)
COLIMITED COLIMITED
ROOMS3101-3105SINGGACOMMERCIALCENTRE ROOMS3101-3105SINGGACOMMERCIALCENTRE
148CONNAUGHT ROADWEST 148CONNAUGHT ROADWEST
WESTERN DISTRICT,Hong Kong WESTERN DISTRICT,Hong Kong
海海海海海海海海海(香香)有有有有 海海海海海海海海海(香香)有有有有
Tel:852-25935622 Tel:852-25935622
____________________________________________________________________________________________________
TERMS OFPAYMENT:60DAYSAMS
TERMS OFSHIPMENT: FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTY Unit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0200 ADP2120ACPZ-R7 4510066029 31,500PC 0.396000 12,474.00
______3_0_1_00_7_8_4_01________________8_0_7_9_70_9_7__________________________________________________________________
Total: 31,500PC 12,474.00
LessDepositPaid: 0.00
NetTotal: 12,474.00
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出出出出出出(編編EN-1102)isattached
INV#90827520,BRAND:ADI,COUNTRYOFORIGIN:SOUTHKOREA
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39Wang KwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 478: 世建发票1.pdf (page 2)
$ws.Cells.Item(478, 1).Value = @'
世建发票1.pdf
'@
$ws.Cells.Item(478, 2).Value = 2
$ws.Cells.Item(478, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90827622",
    "Invoice Date": "28 Aug 2023",
    "Currency": "USD",
    "Amount": 7155.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(478, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90827622",
    "Invoice Date": "28 Aug 2023",
    "Currency": "USD",
    "Amount": 7155.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(478, 5).Value = @'
[]
'@
$ws.Cells.Item(478, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世世世世(香香)有有有有 InvoiceNo. :90827622
Currency :USD
IssuedBy :VickyYan
Salesman :TONYGUO
Date :28Aug2023
Customer Code :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK)
COLIMITED COLIMITED
ROOMS3101-3105SINGGACOMMERCIALCENTRE ROOMS3101-3105SINGGACOMMERCIALCENTRE
148CONNAUGHT ROADWEST 148CONNAUGHT ROADWEST
WESTERN DISTRICT,Hong Kong WESTERN DISTRICT,Hong Kong
海海海海海海海海海(香香)有有有有 海海海海海海海海海(香香)有有有有
Tel:852-25935622 Tel:852-25935622
____________________________________________________________________________________________________
TERMS OFPAYMENT:60DAYSAMS
TERMS OFSHIPMENT: FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTY Unit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0010 ADUCM420BCBZ-RL7 4510079885 1,500PC 4.770000 7,155.00
______3_0_1_11_0_6_1_01________________8_0_7_9_71_9_4__________________________________________________________________
Total: 1,500PC 7,155.00
LessDepositPaid: 0.00
NetTotal: 7,155.00
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出出出出出出(編編EN-1102)isattached
INV#90827622,BRAND:ADI,COUNTRYOFORIGIN: SINGAPORE
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39Wang KwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 479: 世建发票1.pdf (page 3)
$ws.Cells.Item(479, 1).Value = @'
世建发票1.pdf
'@
$ws.Cells.Item(479, 2).Value = 3
$ws.Cells.Item(479, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90827668",
    "Invoice Date": "28 Aug 2023",
    "Currency": "USD",
    "Amount": 7155.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(479, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90827668",
    "Invoice Date": "28 Aug 2023",
    "Currency": "USD",
    "Amount": 7155.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(479, 5).Value = @'
[]
'@
$ws.Cells.Item(479, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世世世世(香香)有有有有 InvoiceNo. :90827668
Currency :USD
IssuedBy :VickyYan
Salesman :TONYGUO
Date :28Aug2023
Customer Code :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK)
COLIMITED COLIMITED
ROOMS3101-3105SINGGACOMMERCIALCENTRE ROOMS3101-3105SINGGACOMMERCIALCENTRE
148CONNAUGHT ROADWEST 148CONNAUGHT ROADWEST
WESTERN DISTRICT,Hong Kong WESTERN DISTRICT,Hong Kong
海海海海海海海海海(香香)有有有有 海海海海海海海海海(香香)有有有有
Tel:852-25935622 Tel:852-25935622
____________________________________________________________________________________________________
TERMS OFPAYMENT:60DAYSAMS
TERMS OFSHIPMENT: FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTY Unit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0010 ADUCM420BCBZ-RL7 4510080078 1,500PC 4.770000 7,155.00
______3_0_1_11_0_6_1_01________________8_0_7_9_71_9_2__________________________________________________________________
Total: 1,500PC 7,155.00
LessDepositPaid: 0.00
NetTotal: 7,155.00
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出出出出出出(編編EN-1102)isattached
INV#90827668,BRAND:ADI,COUNTRYOFORIGIN:SINGAPORE
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39Wang KwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 480: 世建发票1.pdf (page 4)
$ws.Cells.Item(480, 1).Value = @'
世建发票1.pdf
'@
$ws.Cells.Item(480, 2).Value = 4
$ws.Cells.Item(480, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90827701",
    "Invoice Date": "28 Aug 2023",
    "Currency": "USD",
    "Amount": 23836.5,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(480, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90827701",
    "Invoice Date": "28 Aug 2023",
    "Currency": "USD",
    "Amount": 23836.5,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(480, 5).Value = @'
[]
'@
$ws.Cells.Item(480, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世世世世(香香)有有有有 InvoiceNo. :90827701
Currency :USD
IssuedBy :VickyYan
Salesman :TONYGUO
Date :28Aug2023
Customer Code :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK)
COLIMITED COLIMITED
ROOMS3101-3105SINGGACOMMERCIALCENTRE ROOMS3101-3105SINGGACOMMERCIALCENTRE
148CONNAUGHT ROADWEST 148CONNAUGHT ROADWEST
WESTERN DISTRICT,Hong Kong WESTERN DISTRICT,Hong Kong
海海海海海海海海海(香香)有有有有 海海海海海海海海海(香香)有有有有
Tel:852-25935622 Tel:852-25935622
____________________________________________________________________________________________________
TERMS OFPAYMENT:60DAYSAMS
TERMS OFSHIPMENT: FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTY Unit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0010 ADUCM420BCBZ-RL7 4510077212 1,500PC 5.297000 7,945.50
______3_0_1_11_0_6_1_01________________8_0_7_9_72_6_6__________________________________________________________________
0070 ADUCM420BCBZ-RL7 4510077212 1,500PC 5.297000 7,945.50
______3_0_1_11_0_6_1_01________________8_0_7_9_72_6_6__________________________________________________________________
0100 ADUCM420BCBZ-RL7 4510077212 1,500PC 5.297000 7,945.50
______3_0_1_11_0_6_1_01________________8_0_7_9_72_6_6__________________________________________________________________
Total: 4,500PC 23,836.50
LessDepositPaid: 0.00
NetTotal: 23,836.50
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出出出出出出(編編EN-1102)isattached
INV#90827701,BRAND:ADI,COUNTRYOFORIGIN:SINGAPORE
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39Wang KwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 481: 世建发票7.pdf (page 1)
$ws.Cells.Item(481, 1).Value = @'
世建发票7.pdf
'@
$ws.Cells.Item(481, 2).Value = 1
$ws.Cells.Item(481, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90817165",
    "Invoice Date": "03Jul2023",
    "Currency": "USD",
    "Amount": 31782.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(481, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90817165",
    "Invoice Date": "03Jul2023",
    "Currency": "USD",
    "Amount": 31782.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(481, 5).Value = @'
[]
'@
$ws.Cells.Item(481, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世世世世(香香)有有有有 InvoiceNo. :90817165
Currency :USD
IssuedBy :VickyYan
Salesman :TONYGUO
Date :03Jul2023
Customer Code :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK)
COLIMITED COLIMITED
ROOMS3101-3105SINGGACOMMERCIALCENTRE ROOMS3101-3105SINGGACOMMERCIALCENTRE
148CONNAUGHT ROADWEST 148CONNAUGHT ROADWEST
WESTERN DISTRICT,Hong Kong WESTERN DISTRICT,Hong Kong
海海海海海海海海海(香香)有有有有 海海海海海海海海海(香香)有有有有
Tel:852-25935622 Tel:852-25935622
____________________________________________________________________________________________________
TERMS OFPAYMENT:60DAYSAMS
TERMS OFSHIPMENT: FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTY Unit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0030 ADUCM420BCBZ-RL7 4510077212 1,500PC 5.297000 7,945.50
______3_0_1_11_0_6_1_01________________8_0_7_8_69_5_8__________________________________________________________________
0040 ADUCM420BCBZ-RL7 4510077212 1,500PC 5.297000 7,945.50
______3_0_1_11_0_6_1_01________________8_0_7_8_69_5_8__________________________________________________________________
0050 ADUCM420BCBZ-RL7 4510077212 1,500PC 5.297000 7,945.50
______3_0_1_11_0_6_1_01________________8_0_7_8_69_5_8__________________________________________________________________
0060 ADUCM420BCBZ-RL7 4510077212 1,500PC 5.297000 7,945.50
______3_0_1_11_0_6_1_01________________8_0_7_8_69_5_8__________________________________________________________________
Total: 6,000PC 31,782.00
LessDepositPaid: 0.00
NetTotal: 31,782.00
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出出出出出出(編編EN-1102)isattached
INV#90817165,BRAND:ADI,COUNTRYOFORIGIN: SINGAPORE
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39Wang KwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 482: 世建发票3.pdf (page 1)
$ws.Cells.Item(482, 1).Value = @'
世建发票3.pdf
'@
$ws.Cells.Item(482, 2).Value = 1
$ws.Cells.Item(482, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824620",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 4750.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(482, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824620",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 4750.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(482, 5).Value = @'
[]
'@
$ws.Cells.Item(482, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世健系統(香港)有限公司 InvoiceNo. :90824620
Currency :USD
IssuedBy :ShirleyHu
Salesman :TONYGUO
Date :09Aug2023
CustomerCode :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIESCO
COLIMITED LTD
ROOMS3101-3105SINGGACOMMERCIALCENTRE NO218QIANWANROAD
148CONNAUGHTROADWEST QINGDAOECONOMICANDTECHNOLOGICAL
WESTERNDISTRICT,HongKong DEVELOPMENTZONE
海信寬帶多媒體技術(香港)有限公司 QINGDAO
Tel:852-25935622 QINGDAO,China
PostalCode:266071
青島海信寬帶多媒體技術有限公司
青島經濟技術開發區前灣路218號
Tel:86-532-86016016Fax:86-532-86016007
____________________________________________________________________________________________________
TERMSOFPAYMENT:60DAYSAMS
TERMSOFSHIPMENT:FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTYUnit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0370 ADUC7020BCPZ62I-RL 4570065409 2,500PC 1.900000 4,750.00
______3_0_1_00_0_3_4_01________________8_0_7_9_42_9_1__________________________________________________________________
Total: 2,500PC 4,750.00
LessDepositPaid: 0.00
NetTotal: 4,750.00
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出口事宜須知(編號EN-1102)isattached
REF#AL0000036780,INV#90824620,BRAND:ADI,COUNTRYOFORIGIN:CHINA
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39WangKwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 483: 世建发票3.pdf (page 2)
$ws.Cells.Item(483, 1).Value = @'
世建发票3.pdf
'@
$ws.Cells.Item(483, 2).Value = 2
$ws.Cells.Item(483, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824626",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 4338.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(483, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824626",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 4338.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(483, 5).Value = @'
[]
'@
$ws.Cells.Item(483, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世健系統(香港)有限公司 InvoiceNo. :90824626
Currency :USD
IssuedBy :ShirleyHu
Salesman :TONYGUO
Date :09Aug2023
CustomerCode :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIESCO
COLIMITED LTD
ROOMS3101-3105SINGGACOMMERCIALCENTRE NO218QIANWANROAD
148CONNAUGHTROADWEST QINGDAOECONOMICANDTECHNOLOGICAL
WESTERNDISTRICT,HongKong DEVELOPMENTZONE
海信寬帶多媒體技術(香港)有限公司 QINGDAO
Tel:852-25935622 QINGDAO,China
PostalCode:266071
青島海信寬帶多媒體技術有限公司
青島經濟技術開發區前灣路218號
Tel:86-532-86016016Fax:86-532-86016007
____________________________________________________________________________________________________
TERMSOFPAYMENT:60DAYSAMS
TERMSOFSHIPMENT:FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTYUnit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0010 ADA4096-4ACPZ-R7 4570067453 1,500PC 2.892000 4,338.00
______3_0_1_10_3_0_5_01________________8_0_7_9_42_9_7__________________________________________________________________
Total: 1,500PC 4,338.00
LessDepositPaid: 0.00
NetTotal: 4,338.00
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出口事宜須知(編號EN-1102)isattached
REF#AL0000037196,INV#90824626,BRAND:ADI,COUNTRYOFORIGIN:KOREA
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39WangKwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 484: 世建发票3.pdf (page 3)
$ws.Cells.Item(484, 1).Value = @'
世建发票3.pdf
'@
$ws.Cells.Item(484, 2).Value = 3
$ws.Cells.Item(484, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824632",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 30270.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(484, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824632",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 30270.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(484, 5).Value = @'
[]
'@
$ws.Cells.Item(484, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世健系統(香港)有限公司 InvoiceNo. :90824632
Currency :USD
IssuedBy :ShirleyHu
Salesman :TONYGUO
Date :09Aug2023
CustomerCode :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIESCO
COLIMITED LTD
ROOMS3101-3105SINGGACOMMERCIALCENTRE NO218QIANWANROAD
148CONNAUGHTROADWEST QINGDAOECONOMICANDTECHNOLOGICAL
WESTERNDISTRICT,HongKong DEVELOPMENTZONE
海信寬帶多媒體技術(香港)有限公司 QINGDAO
Tel:852-25935622 QINGDAO,China
PostalCode:266071
青島海信寬帶多媒體技術有限公司
青島經濟技術開發區前灣路218號
Tel:86-532-86016016Fax:86-532-86016007
____________________________________________________________________________________________________
TERMSOFPAYMENT:60DAYSAMS
TERMSOFSHIPMENT:FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTYUnit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0040 AD5679RBCPZ-2-RL7 4570067453 1,500PC 20.180000 30,270.00
______3_0_1_11_3_0_4_01________________8_0_7_9_43_0_2__________________________________________________________________
Total: 1,500PC 30,270.00
LessDepositPaid: 0.00
NetTotal: 30,270.00
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出口事宜須知(編號EN-1102)isattached
REF#AL0000037196,INV#90824632,BRAND:ADI,COUNTRYOFORIGIN:SOUTH
KOREA
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39WangKwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 485: 世建发票3.pdf (page 4)
$ws.Cells.Item(485, 1).Value = @'
世建发票3.pdf
'@
$ws.Cells.Item(485, 2).Value = 4
$ws.Cells.Item(485, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824629",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 5576.69,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(485, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824629",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 5576.69,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(485, 5).Value = @'
[]
'@
$ws.Cells.Item(485, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世健系統(香港)有限公司 InvoiceNo. :90824629
Currency :USD
IssuedBy :ShirleyHu
Salesman :TONYGUO
Date :09Aug2023
CustomerCode :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIESCO
COLIMITED LTD
ROOMS3101-3105SINGGACOMMERCIALCENTRE NO218QIANWANROAD
148CONNAUGHTROADWEST QINGDAOECONOMICANDTECHNOLOGICAL
WESTERNDISTRICT,HongKong DEVELOPMENTZONE
海信寬帶多媒體技術(香港)有限公司 QINGDAO
Tel:852-25935622 QINGDAO,China
PostalCode:266071
青島海信寬帶多媒體技術有限公司
青島經濟技術開發區前灣路218號
Tel:86-532-86016016Fax:86-532-86016007
____________________________________________________________________________________________________
TERMSOFPAYMENT:60DAYSAMS
TERMSOFSHIPMENT:FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTYUnit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0030 LTC3311SEV#PBF 4570066240 1,872PC 2.979000 5,576.69
______3_0_1_11_1_5_8_01________________8_0_7_9_43_0_0__________________________________________________________________
Total: 1,872PC 5,576.69
LessDepositPaid: 0.00
NetTotal: 5,576.69
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出口事宜須知(編號EN-1102)isattached
REF#AL0000036955,INV#90824629,BRAND:ADI,COUNTRYOFORIGIN:SOUTH
KOREA
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39WangKwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 486: 世建发票3.pdf (page 5)
$ws.Cells.Item(486, 1).Value = @'
世建发票3.pdf
'@
$ws.Cells.Item(486, 2).Value = 5
$ws.Cells.Item(486, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824628",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 873.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(486, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824628",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 873.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(486, 5).Value = @'
[]
'@
$ws.Cells.Item(486, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世健系統(香港)有限公司 InvoiceNo. :90824628
Currency :USD
IssuedBy :ShirleyHu
Salesman :TONYGUO
Date :09Aug2023
CustomerCode :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIESCO
COLIMITED LTD
ROOMS3101-3105SINGGACOMMERCIALCENTRE NO218QIANWANROAD
148CONNAUGHTROADWEST QINGDAOECONOMICANDTECHNOLOGICAL
WESTERNDISTRICT,HongKong DEVELOPMENTZONE
海信寬帶多媒體技術(香港)有限公司 QINGDAO
Tel:852-25935622 QINGDAO,China
PostalCode:266071
青島海信寬帶多媒體技術有限公司
青島經濟技術開發區前灣路218號
Tel:86-532-86016016Fax:86-532-86016007
____________________________________________________________________________________________________
TERMSOFPAYMENT:60DAYSAMS
TERMSOFSHIPMENT:FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTYUnit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0050 LTC3308AEV#TRMPBF 4570066645 500PC 1.746000 873.00
______3_0_1_11_0_8_1_01________________8_0_7_9_42_9_9__________________________________________________________________
Total: 500PC 873.00
LessDepositPaid: 0.00
NetTotal: 873.00
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出口事宜須知(編號EN-1102)isattached
REF#AL0000037071,INV#90824628,BRAND:ADI,COUNTRYOFORIGIN:SOUTH
KOREA
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39WangKwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 487: 世建发票3.pdf (page 6)
$ws.Cells.Item(487, 1).Value = @'
世建发票3.pdf
'@
$ws.Cells.Item(487, 2).Value = 6
$ws.Cells.Item(487, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824625",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 5104.5,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(487, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "90824625",
    "Invoice Date": "09 Aug 2023",
    "Currency": "USD",
    "Amount": 5104.5,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO LIMITED",
    "From": "Excelpoint Systems (H.K.) Limited",
    "page": "1/1"
}
'@
$ws.Cells.Item(487, 5).Value = @'
[]
'@
$ws.Cells.Item(487, 6).Value = @'
ExcelpointSystems(H.K.)Limited Page :1/ 1
世健系統(香港)有限公司 InvoiceNo. :90824625
Currency :USD
IssuedBy :ShirleyHu
Salesman :TONYGUO
Date :09Aug2023
CustomerCode :6000004601
____________________________________________________________________________________________________
INVOICE
____________________________________________________________________________________________________
BILLTO: SOLDTO:
HISENSEBROADBANDMULTIMEDIATECHNOLOGIES(HK) HISENSEBROADBANDMULTIMEDIATECHNOLOGIESCO
COLIMITED LTD
ROOMS3101-3105SINGGACOMMERCIALCENTRE NO218QIANWANROAD
148CONNAUGHTROADWEST QINGDAOECONOMICANDTECHNOLOGICAL
WESTERNDISTRICT,HongKong DEVELOPMENTZONE
海信寬帶多媒體技術(香港)有限公司 QINGDAO
Tel:852-25935622 QINGDAO,China
PostalCode:266071
青島海信寬帶多媒體技術有限公司
青島經濟技術開發區前灣路218號
Tel:86-532-86016016Fax:86-532-86016007
____________________________________________________________________________________________________
TERMSOFPAYMENT:60DAYSAMS
TERMSOFSHIPMENT:FCAHONGKONG
_____________________________________________________________________________________________________
Item# Productcode PONo QTYUnit Unitprice Amount
ExcelPointP/N DNNo
CustomerP/N
_____________________________________________________________________________________________________
0010 ADN8834ACBZ-R7 4570065810 1,500PC 3.403000 5,104.50
______3_0_1_01_1_8_3_01________________8_0_7_9_42_9_4__________________________________________________________________
Total: 1,500PC 5,104.50
LessDepositPaid: 0.00
NetTotal: 5,104.50
======================================================================================
____________________________________________________________________________________________________
ThisInvoiceisacomputerprintoutandnosignatureisrequired
SHIPPINGMARKS:
Remarks:
出口事宜須知(編號EN-1102)isattached
REF#AL0000036866,INV#90824625,BRAND:ADI,COUNTRYOFORIGIN:SINGAPORE
____________________________________________________________________________________________________
Address:3108,31/F,SkylineTower,39WangKwongRoad,,HongKong
Phone:(852)2503-2212FAX:(852)2503-1558
'@

# Row 488: WIN TREND XSCKD-SD-230818-09333.pdf (page 1)
$ws.Cells.Item(488, 1).Value = @'
WIN TREND XSCKD-SD-230818-09333.pdf
'@
$ws.Cells.Item(488, 2).Value = 1
$ws.Cells.Item(488, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "XSCKD-SD-230818-09333",
    "Invoice Date": "2023-08-18",
    "Currency": "USD",
    "Amount": 1078.0,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co., Ltd",
    "From": "Win Trend Electronic Technology Limited"
}
'@
$ws.Cells.Item(488, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "XSCKD-SD-230818-09333",
    "Invoice Date": "2023-08-18",
    "Currency": "USD",
    "Amount": 1078.0,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co., Ltd",
    "From": "Win Trend Electronic Technology Limited"
}
'@
$ws.Cells.Item(488, 5).Value = @'
[]
'@
$ws.Cells.Item(488, 6).Value = @'
INVOICE
Win Trend Electronic Technology Limited.
6/F, EVERWIN CENTRE, NO.72 HUNG TO ROAD, KWUN TONG, KOWLOON, HONG KONG
Invoice no.: XSCKD-SD-230818-09333
Tel :Sally Siu 9623 7440 / Angie Kwan 9408 3888
Bill to: Hisense Broadband Multimedia Technologies (HK) Co.,Ltd Date: 2023-08-18
Attn. To: Shipping Terms
Tel.: +86(0)532-80879676 / 18366214997 Payment Terms:月结60天
Room 3101-3105, Singga Commercial Centre, 148 Connaught
Add.: Sales: 孙志鹏
Rd. West, Sai Ying Pun, HK
QTY U/P AMT
ITEM CSTM P/O NO. BRAND CSTM P/N DESCRIPTION
PCS USD USD
1 4570072564 CREDO/默升 3011105301 CFD50502-A0-ABXN 98 11.00 1078.00
合计： 98 1078.00
Beneficiary:Win Trend Electronic Technology Ltd.
Bank Name:The Hong Kong and Shanghai Banking Corporation Ltd.
Add.: No.1 Queen's Road Central, Hong Kong
A/C No.: 848 533 097 838
Swift Code:HSBCHKHHHKH
送貨地址: 青旅思捷物流有限公司CYTS-SPIRIT LOGISTICS LIMITED
香港荃灣橫窩仔街43-57號永得利中心六樓6/F, Ever Gain Centre, 43-57 Wang Wo Tsai Street, Tsuen Wan, N.T., Hong Kong
TEL: 31522431 联系人： Mr. Key (阿仁)
入仓单号： ZGD-506019 备注： 江门
E. &O. E.
CONFIRM & RECEIVED BY:
'@

# Row 489: WIN TREND XSCKD-SD-230818-09333.pdf (page 2)
$ws.Cells.Item(489, 1).Value = @'
WIN TREND XSCKD-SD-230818-09333.pdf
'@
$ws.Cells.Item(489, 2).Value = 2
$ws.Cells.Item(489, 3).Value = @'
{
    "Doc Type": "非发票：可能是装货单、waybill或其他"
}
'@
$ws.Cells.Item(489, 4).Value = @'
{
    "Doc Type": "非发票：可能是装货单、waybill或其他"
}
'@
$ws.Cells.Item(489, 5).Value = @'
[]
'@
$ws.Cells.Item(489, 6).Value = @'
PACKING LIST
Win Trend Electronic Technology Limited
6/F, EVERWIN CENTRE, NO.72 HUNG TO ROAD, KWUN TONG, KOWLOON, HONG KONG
Tel :Sally Siu 9623 7440 / Angie Kwan 9408 3888
Bill to: Date: 2023-08-18
Hisense Broadband Multimedia Technologies (HK) Co.,Ltd Invoice No. :XSCKD-SD-230818-09333
Ship to: TEL: 31522431
青旅思捷物流有限公司
CYTS-SPIRIT LOGISTICS LIMITED
Contact: Mr. Key (阿仁)
香港荃灣橫窩仔街43-57號永得利中心六樓
6/F, Ever Gain Centre, 43-57 Wang Wo Tsai Street, Tsuen Wan, N.T., Hong
Kong
N.W. G.W.
QTY SIZE
CTN No. PART No. CSTM PART No. PO NO. BRAND COO /CTN /CTN Remark
(PCS) ( CM )
KG KG
CREDO/默
1 CFD50502-A0-ABXN 3011105301 4570072564 KR 98 0.9 1.3 40X40X15 江门
升
合计： 98
for and on behalf of for and on behalf of
WIN TREND Electronic Technology Limited.
Authorized Signatures Confirmed & Received by
'@

# Row 490: WIN TREND XSCKD-SD-230818-09341.pdf (page 1)
$ws.Cells.Item(490, 1).Value = @'
WIN TREND XSCKD-SD-230818-09341.pdf
'@
$ws.Cells.Item(490, 2).Value = 1
$ws.Cells.Item(490, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "XSCKD-SD-230818-09341",
    "Invoice Date": "2023-08-18",
    "Currency": "USD",
    "Amount": 0.0,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co.,Ltd",
    "From": "Win Trend Electronic Technology Limited"
}
'@
$ws.Cells.Item(490, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "XSCKD-SD-230818-09341",
    "Invoice Date": "2023-08-18",
    "Currency": "USD",
    "Amount": 0.0,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co.,Ltd",
    "From": "Win Trend Electronic Technology Limited"
}
'@
$ws.Cells.Item(490, 5).Value = @'
[]
'@
$ws.Cells.Item(490, 6).Value = @'
INVOICE
Win Trend Electronic Technology Limited.
6/F, EVERWIN CENTRE, NO.72 HUNG TO ROAD, KWUN TONG, KOWLOON, HONG KONG
Invoice no.: XSCKD-SD-230818-09341
Tel :Sally Siu 9623 7440 / Angie Kwan 9408 3888
Bill to: Hisense Broadband Multimedia Technologies (HK) Co.,Ltd Date: 2023-08-18
Attn. To: Shipping Terms:
Tel.: +86(0)532-80879676 / 18366214997 Payment Terms: 月结60天
Room 3101-3105, Singga Commercial Centre, 148 Connaught
Add.: Sales: 孙志鹏
Rd. West, Sai Ying Pun, HK
QTY U/P AMT
ITEM CSTM P/O NO. BRAND CSTM P/N DESCRIPTION
PCS USD USD
1 5010001959 CREDO/默升 3011105301 CFD50502-A0-ABXN 120 0.00 0.00
合计： 120 0.00
Beneficiary:Win Trend Electronic Technology Ltd.
Bank Name:The Hong Kong and Shanghai Banking Corporation Ltd.
Add.: No.1 Queen's Road Central, Hong Kong
A/C No.: 848 533 097 838
Swift Code:HSBCHKHHHKH
送貨地址: 青旅思捷物流有限公司CYTS-SPIRIT LOGISTICS LIMITED
香港荃灣橫窩仔街43-57號永得利中心六樓6/F, Ever Gain Centre, 43-57 Wang Wo Tsai Street, Tsuen Wan, N.T., Hong Kong
TEL: 31522431 联系人： Mr. Key (阿仁)
ZGD-506021 江门
入仓单号： 备注：
E. &O. E.
CONFIRM & RECEIVED BY:
'@

# Row 491: WIN TREND XSCKD-SD-230818-09341.pdf (page 2)
$ws.Cells.Item(491, 1).Value = @'
WIN TREND XSCKD-SD-230818-09341.pdf
'@
$ws.Cells.Item(491, 2).Value = 2
$ws.Cells.Item(491, 3).Value = @'
{
    "Doc Type": "非发票：可能是装货单、waybill或其他"
}
'@
$ws.Cells.Item(491, 4).Value = @'
{
    "Doc Type": "非发票：可能是装货单、waybill或其他"
}
'@
$ws.Cells.Item(491, 5).Value = @'
[]
'@
$ws.Cells.Item(491, 6).Value = @'
PACKING LIST
Win Trend Electronic Technology Limited
6/F, EVERWIN CENTRE, NO.72 HUNG TO ROAD, KWUN TONG, KOWLOON, HONG KONG
Tel :Sally Siu 9623 7440 / Angie Kwan 9408 3888
Bill to: Date: 2023-08-18
Hisense Broadband Multimedia Technologies (HK) Co.,Ltd Invoice No. : XSCKD-SD-230818-09341
Ship to: TEL: 31522431
青旅思捷物流有限公司
CYTS-SPIRIT LOGISTICS LIMITED
Contact: Mr. Key (阿仁)
香港荃灣橫窩仔街43-57號永得利中心六樓
6/F, Ever Gain Centre, 43-57 Wang Wo Tsai Street, Tsuen Wan, N.T., Hong
Kong
N.W. G.W.
QTY SIZE
CTN No. PART No. CSTM PART No. PO NO. BRAND COO /CTN /CTN Remark
(PCS) ( CM )
KG KG
CREDO/默
1 CFD50502-A0-ABXN 3011105301 5010001959 TW,CN 120 0.7 1.4 40X40X15 江门
升
合计： 1 120 0.7 1.4
for and on behalf of for and on behalf of
WIN TREND Electronic Technology Limited.
Authorized Signatures Confirmed & Received by
'@

# Row 492: WIN TREND XSCKD-SD-230907-09490.pdf (page 1)
$ws.Cells.Item(492, 1).Value = @'
WIN TREND XSCKD-SD-230907-09490.pdf
'@
$ws.Cells.Item(492, 2).Value = 1
$ws.Cells.Item(492, 3).Value = @'
{
    "Doc Type": "其他",
    "Invoice No.": "",
    "Invoice Date": "2023-09-07",
    "Currency": "",
    "Amount": 0,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co.,Ltd",
    "From": "WIN TREND Electronic Technology Limited"
}
'@
$ws.Cells.Item(492, 4).Value = @'
{
    "Doc Type": "其他",
    "Invoice No.": "",
    "Invoice Date": "2023-09-07",
    "Currency": "",
    "Amount": 0,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co.,Ltd",
    "From": "WIN TREND Electronic Technology Limited"
}
'@
$ws.Cells.Item(492, 5).Value = @'
[]
'@
$ws.Cells.Item(492, 6).Value = @'
PACKING
G LIST
胜达科
WIN TREND TECHNOLOGY
Win Trend ElectronicTechnology Limited
6/F,EVERWINCENTRE,NO.72HUNGTOROAD,KWUNTONG,KOWLOON,HONGKONG
Tel:Sally Siu 96237440/Angie Kwan 94083888
Bill to:
Date:
2023-09-07
Hisense Broadband Multimedia Technologies (HK) Co.,Ltd
Invoice No.: XSCKD-SD-230907-09490
Ship to:
TEL:
31522431
青旅思捷物流有限公司
CYTS-SPIRIT LOGISTICS LIMITED
Contact:
Mr.Key（阿仁）
香港荃灣横商仔街43-57號永得利中心六楼
6/F. Ever Gain Centre. 43-57 Wang Wo Tsai Street. Tsuen Wan. N.T..
N.W.
G.W.
BRAND
QTY
SIZE
CTN No.
PART No.
CSTM PART No.
PO NO.
COO
/CTN
/CTN
(PCS)
(CM）
KG
KG
1
MP2161GJ-Z
3010074101
/
MPS/芯源
CHINA
6000
0.5
0.8
24X23X14
合计：
1
0.5
0.8
for and on behalf of
for and on behalf of
WIN TREND Electronic Technology Limited.
Authorized Signatures
Confirmed & Received by
'@

# Row 493: WIN TREND XSCKD-SD-230907-09490.pdf (page 2)
$ws.Cells.Item(493, 1).Value = @'
WIN TREND XSCKD-SD-230907-09490.pdf
'@
$ws.Cells.Item(493, 2).Value = 2
$ws.Cells.Item(493, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "XSCKD-SD-230907-09490",
    "Invoice Date": "2023-09-07",
    "Currency": "USD",
    "Amount": 559.8,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co.,Ltd",
    "From": "Win Trend Electronic Technology Limited"
}
'@
$ws.Cells.Item(493, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "XSCKD-SD-230907-09490",
    "Invoice Date": "2023-09-07",
    "Currency": "USD",
    "Amount": 559.8,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co.,Ltd",
    "From": "Win Trend Electronic Technology Limited"
}
'@
$ws.Cells.Item(493, 5).Value = @'
[]
'@
$ws.Cells.Item(493, 6).Value = @'
INVOICE
WIN TREMD FECHBNO)IOGY
Win Trend Electronic Technology Limited.
6/F, EVERWIN CENTRE, NO.72 HUNG TO ROAD, KWUN TONG, KOWLOON, HONG KONG
Invoice no.:
XSCKD-SD-230907-09490
Tel :Sally Siu 9623 7440 / Angie Kwan 9408 3888
Bill to:
Hisense Broadband Multimedia Technologies (HK) Co.,Ltd
Date:
2023-09-07
Attn. To:
Shipping
Tel. :
+86 (0)532-80879676 / 18366214997
Payment
月结60天
Room 3101-3105, Singga Commercial Centre,148 Connaught Rd.
Add. :
Sales:
孙志鹏
West, Sai Ying Pun, HK
QTY
d/n
AMT
ITEM
CSTM P/0 NO.
BRAND
CSTM P/N
DESCRIPTION
PCS
USD
USD
1
MPS/芯源
3010074101
MP2161GJ-Z
6000
0.0933
559.80
Beneficiary:Win Trend Electronic Technology Ltd.
Bank Name:The Hong Kong and Shanghai Banking Corporation Ltd.
Add.: No.1 Queen' s Road Central, Hong Kong
A/C No.:848533097838
Swift Code:HSBCHKHHHKH
送货地址：
青旅思捷物流有限公司CYTS-SPIRIT LOGISTICSLIMITED
香港基灣横仔街43-57號永得利中心六樓6/F，Ever Gain Centre，43-57Wang Wo Tsai Street，Tsuen Wan，N.T.，Hong Kong
TEL:
31522431
联系人：
Mr.Key (阿仁)
入仓单号：
NA
备注：
青岛
E.&O.E.
CONFIRM & RECEIVED BY:
'@

# Row 494: XSCKD-SD-230818-09332 INV-青岛.pdf (page 1)
$ws.Cells.Item(494, 1).Value = @'
XSCKD-SD-230818-09332 INV-青岛.pdf
'@
$ws.Cells.Item(494, 2).Value = 1
$ws.Cells.Item(494, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "XSCKD-SD-230818-09332",
    "Invoice Date": "2023-08-18",
    "Currency": "USD",
    "Amount": 5170.0,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co., Ltd",
    "From": "Win Trend Electronic Technology Limited"
}
'@
$ws.Cells.Item(494, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "XSCKD-SD-230818-09332",
    "Invoice Date": "2023-08-18",
    "Currency": "USD",
    "Amount": 5170.0,
    "Bill To": "Hisense Broadband Multimedia Technologies (HK) Co., Ltd",
    "From": "Win Trend Electronic Technology Limited"
}
'@
$ws.Cells.Item(494, 5).Value = @'
[]
'@
$ws.Cells.Item(494, 6).Value = @'
INVOICE
Win Trend Electronic Technology Limited.
6/F, EVERWIN CENTRE, NO.72 HUNG TO ROAD, KWUN TONG, KOWLOON, HONG KONG
Invoice no.: XSCKD-SD-230818-09332
Tel :Sally Siu 9623 7440 / Angie Kwan 9408 3888
Bill to: Hisense Broadband Multimedia Technologies (HK) Co.,Ltd Date: 2023-08-18
Attn. To: Shipping Terms:
Tel.: +86(0)532-80879676 / 18366214997 Payment Terms: 月结60天
Room 3101-3105, Singga Commercial Centre, 148 Connaught
Add.: Sales: 孙志鹏
Rd. West, Sai Ying Pun, HK
QTY U/P AMT
ITEM CSTM P/O NO. BRAND CSTM P/N DESCRIPTION
PCS USD USD
1 4570072564 CREDO/默升 3011105301 CFD50502-A0-ABXN 470 11.00 5170.00
合计： 470 5170.00
Beneficiary:Win Trend Electronic Technology Ltd.
Bank Name:The Hong Kong and Shanghai Banking Corporation Ltd.
Add.: No.1 Queen's Road Central, Hong Kong
A/C No.: 848 533 097 838
Swift Code:HSBCHKHHHKH
送貨地址: 青旅思捷物流有限公司CYTS-SPIRIT LOGISTICS LIMITED
香港荃灣橫窩仔街43-57號永得利中心六樓6/F, Ever Gain Centre, 43-57 Wang Wo Tsai Street, Tsuen Wan, N.T., Hong Kong
TEL: 31522431 联系人： Mr. Key (阿仁)
入仓单号： NA 备注： 青岛
E. &O. E.
CONFIRM & RECEIVED BY:
'@

# Row 495: HK00213203.pdf (page 1)
$ws.Cells.Item(495, 1).Value = @'
HK00213203.pdf
'@
$ws.Cells.Item(495, 2).Value = 1
$ws.Cells.Item(495, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "HK00213203",
    "Invoice Date": "2023/4/13",
    "Currency": "USD",
    "Amount": 27000.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO., LTD.",
    "From": "HITACHI HIGH-TECH HONG KONG LTD",
    "page": "1/1"
}
'@
$ws.Cells.Item(495, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "HK00213203",
    "Invoice Date": "2023/4/13",
    "Currency": "USD",
    "Amount": 27000.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO., LTD.",
    "From": "HITACHI HIGH-TECH HONG KONG LTD",
    "page": "1/1"
}
'@
$ws.Cells.Item(495, 5).Value = @'
[]
'@
$ws.Cells.Item(495, 6).Value = @'
Hitachi High-Tech Hong Kong Ltd.
8/F, Building 20E, Phase 3, Hong Kong Science Park,
Pak Shek Kok, N.T., Hong Kong
TEL: 852-2737-4700
FAX: 852-2370-7255
Invoice
Bill To: PAGE 1/1
Purchasing Departmet Invoice No.: HK00213203
HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO.,LTD. Invoice Date : 2023/4/13
ROOMS 3010-3105 SINGGA COMMERCIAL CENTRE Inco Terms : FCA JAPAN
148 CONNAUGHT ROAD WEST HONG KONG Terms Of Payment : T/T 60 DAYS MONTYLY
Shipped Per : DHL
Ship To : Shipment date: 2023/4/13
Loading Port : TOKYO
CYTS-SPIRIT LOGISTICS LIMITED Discharge Port : HONG KONG
Eternity Warehouse 6/F, Ever Gain Centre, 43-57 Wang Wo Tsai Street, Vendor Code: 5001060
Tsuen Wan, N.T.. HONGKONG
boyim@cyts-spirit.com
Item Cust. PO No. Description Part NO. Quantity Unit Price Amount
(PCS) (USD) (USD)
1 4510070120 4431004801 1-P7175-01-R00 100,000 0.27 27,000.00
Total : 100,000 27,000.00
Packing list
COUNTRY Gross
OF ORIGIN Quantity(PCS) Carton Net Weight(KGS) Weight(KGS) CBM(M3)
JAPAN 4431004801 50,000 2-1 0.10 3.30 0.017204
4431004801 50,000 2-2 0.10 3.30 0.017204
Total: 100,000 2 0.20 6.60 0.034
Remarks:
入仓号: ZGD-505007 (广东宽带多媒体货物)
Remit to : MUFG BANK LTD HONG KONG BRANCH
Account Name : HITACHI HIGH-TECH HONG KONG LTD
Account No : 047-821-85806018194
Swift Code : BOTKHKHH
for Hitachi High-Tech Hong Kong Ltd.
Authorized Signature
'@

# Row 496: HK00214688.pdf (page 1)
$ws.Cells.Item(496, 1).Value = @'
HK00214688.pdf
'@
$ws.Cells.Item(496, 2).Value = 1
$ws.Cells.Item(496, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "HK00214688",
    "Invoice Date": "2023/7/7",
    "Currency": "USD",
    "Amount": 4800.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO., LTD.",
    "From": "HITACHI HIGH-TECH HONG KONG LTD",
    "page": "1/1"
}
'@
$ws.Cells.Item(496, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "HK00214688",
    "Invoice Date": "2023/7/7",
    "Currency": "USD",
    "Amount": 4800.0,
    "Bill To": "HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO., LTD.",
    "From": "HITACHI HIGH-TECH HONG KONG LTD",
    "page": "1/1"
}
'@
$ws.Cells.Item(496, 5).Value = @'
[]
'@
$ws.Cells.Item(496, 6).Value = @'
Hitachi High-Tech Hong Kong Ltd.
8/F, Building 20E, Phase 3, Hong Kong Science Park,
Pak Shek Kok, N.T., Hong Kong
TEL: 852-2737-4700
FAX: 852-2370-7255
Invoice
Bill To: PAGE 1/1
Purchasing Departmet Invoice No.: HK00214688
HISENSE BROADBAND MULTIMEDIA TECHNOLOGIES (HK) CO.,LTD. Invoice Date : 2023/7/7
ROOMS 3010-3105 SINGGA COMMERCIAL CENTRE Inco Terms : FCA JAPAN
148 CONNAUGHT ROAD WEST HONG KONG Terms Of Payment : T/T 60 DAYS MONTYLY
Shipped Per : DHL
Ship To : Shipment date: 2023/7/7
Loading Port : TOKYO
CYTS-SPIRIT LOGISTICS LIMITED Discharge Port : HONG KONG
Eternity Warehouse 6/F, Ever Gain Centre, 43-57 Wang Wo Tsai Street, Vendor Code: 5001060
Tsuen Wan, N.T.. HONGKONG
boyim@cyts-spirit.com
Item Cust. PO No. Description Part NO. Quantity Unit Price Amount
(PCS) (USD) (USD)
1 4570071385 KSHH018XCA 4501003801 1,000 4.80 4,800.00
Total : 1,000 4,800.00
Packing list
COUNTRY Gross
OF ORIGIN Quantity(PCS) Carton Net Weight(KGS) Weight(KGS) CBM(M3)
JAPAN KSHH018XCA 1,000 1-1 0.05 0.70 0.01
Total: 1,000 1 0.05 0.70 0.010
Remarks:
Remit to : MUFG BANK LTD HONG KONG BRANCH
Account Name : HITACHI HIGH-TECH HONG KONG LTD
Account No : 047-821-85806018194
Swift Code : BOTKHKHH
for Hitachi High-Tech Hong Kong Ltd.
Authorized Signature
'@

# Row 497: MISSION 5970023769.pdf (page 1)
$ws.Cells.Item(497, 1).Value = @'
MISSION 5970023769.pdf
'@
$ws.Cells.Item(497, 2).Value = 1
$ws.Cells.Item(497, 3).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "F230726000015",
    "Invoice Date": "2023-07-26",
    "Currency": "",
    "Amount": 0,
    "Bill To": "Hisense broadband multimedia technologies (HK) CO., Ltd",
    "From": "Mission Electronics International Limited"
}
'@
$ws.Cells.Item(497, 4).Value = @'
{
    "Doc Type": "Invoice",
    "Invoice No.": "F230726000015",
    "Invoice Date": "2023-07-26",
    "Currency": "",
    "Amount": 0,
    "Bill To": "Hisense broadband multimedia technologies (HK) CO., Ltd",
    "From": "Mission Electronics International Limited"
}
'@
$ws.Cells.Item(497, 5).Value = @'
[]
'@
$ws.Cells.Item(497, 6).Value = @'
美盛電子國際有限公司
MISSION ELECTRONICSINTERNATIONAL LIMITED
Webite:www.mission-io.com
General line:((00852)27582283Fax:(00852)31051139
美盛
盛Address:香港九龍九龍灣宏開道20號杨耀崧第八工業大厦9樓E室
NOISSIW
Flat E,9/F,Yeung Yiu Chung (No.8) Industrial Building,No.20 Wang Hoi Road, Kowloon Bay,Kowloon
57023272
PACKINGLIST
Customer:Hisense broadband multimedia technologies (HK) CO.,Ltd
P.O.No.:5970023769
InvoiceNo.:F230726000015
PaymentTerm:发货单出货后30天付款(NET30)
InvoiceDate:2023-07-26
ShipVia:
BillTo:Hisense broadbandmultimedia technologies(HK) CO.,Ltd
ShipTo:CYTS-SPIRITLOGISTICS,LTD
Room 3101-3105, Singga Commercial Centre, 148 Conaught R d.
6/F,43-57.Wang Wo Tsai Street, Ever Gain centre ,Tsuen wan,
West, Sai Ying Pun, HK
NT,HONG KONG
Attn:蔺艳莉
Attn:BO YIM
Tel:15684730865
Tel:31522507
数量
重量
尺寸
产地
行号
客户料号Cusr
规格型号
DC
Brand
QTY
Weight
Dimension
Country of
LOT NO.
CTN#
Internalpart No.
Part No.
年份
品牌
(pcs)
(KG)
(cm)
origin
NW:2.3
3011093801
INA216A3YFFT
750
60*24*24
2325
T1
1
3431356CL3
GW:3.7
3011093801
INA216A3YFFT
1500
2318
Il.
3304605CL3
INA216A3YFFT
500
2318
TI
3011093801
3304608CL3
3011093801
INA216A3YFFT
2500
2318
T
3304607CL3
3011093801
INA216A3YFFT
500
2318
T
3304606CL3
3011093801
INA216A3YFFT
250
2318
T
3304609CL3
NW:1.2
3011089901
TPS62480RNCT
6000
42*15*43
2141
TI
2141
2
GW:2.0
3011089901
3000
2150
2150
TPS62480RNCT
12000
NW:2.8
3
3011089901
TPS62480RNCT
40*30*40
2150
T1
2150
GW:4.3
3000
2204
TI
3011089901
TPS62480RNCT
2204
3011089901
TPS62480RNCT
6000
2226
TI
2226
Total:
3
Cartons
36000 PCS
NW: 6.3 KG ; GW: 10.0 KG
Note:
1.Purchase Order Number must be quoted on all documents to facilitate payment.
2.No Cancellation ,No Return ,No Reschedule (except for Electronic Failure).
3.Any complaints against this invoice must be notified to the company within 7 days of data hereof.
4.We are not responsible for any banking charges outside ShenZhen The Payment received in our
account must be the same as the amount we invoiced.
For and on behalf of
Confirmed and received by
Mission ElectronicsInternational Limited
Viki欧阳桃
英盛
请签章
客子国際
有限公司
Customer Chop and Signature
Authorized Signature
'@

$ws.Range("A477:F497").EntireRow.AutoFit()  # clear auto row-height metadata so new rows match existing plain <row> style

